$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date in column C for every data row (2-246)
# from 45192 (2023-09-13) to 45202 (2023-09-23).
for ($r = 2; $r -le 246; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Row 246 gains an explicit row height (15, custom) in the new version.
$ws.Rows.Item(246).RowHeight = 15

# Append the new record as row 247.
$ws.Cells.Item(247, 1).Value = "A 46122-2023"

$ws.Cells.Item(247, 2).Value = 45196
$ws.Cells.Item(247, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(247, 3).Value = 45202
$ws.Cells.Item(247, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(247, 4).Value = "KRONOBERGS LÄN"
$ws.Cells.Item(247, 5).Value = "MARKARYD"

$ws.Cells.Item(247, 7).Value = 1.6
$ws.Cells.Item(247, 8).Value = 0
$ws.Cells.Item(247, 9).Value = 0
$ws.Cells.Item(247, 10).Value = 0
$ws.Cells.Item(247, 11).Value = 0
$ws.Cells.Item(247, 12).Value = 0
$ws.Cells.Item(247, 13).Value = 0
$ws.Cells.Item(247, 14).Value = 0
$ws.Cells.Item(247, 15).Value = 0
$ws.Cells.Item(247, 16).Value = 0
$ws.Cells.Item(247, 17).Value = 0

$ws.Cells.Item(247, 18).Value = ""
$ws.Cells.Item(247, 18).WrapText = $true
